$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date-like text (e.g. "2021-01-09") stored as plain strings,
# not real dates. A leading apostrophe forces text entry; resetting the
# style back to "Normal" afterwards drops the quote-prefix formatting so the
# cell ends up with the same (default) style as its neighbours.

# Row 53: same "2021-01-09" / "10 Jan -- 16 Jan 2021" prediction block as row
# 50, but with updated J/K (Weekly MAE / Weekly MAPE) values from new weather
# data.
$ws.Range("A53").Value = "'2021-01-09"
$ws.Range("A53").Style = "Normal"
$ws.Range("B53").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("C53").Value = 3333.57
$ws.Range("D53").Value = 2118.92
$ws.Range("E53").Value = 1214.65
$ws.Range("F53").Value = "KNN"
$ws.Range("J53").Value = 1173.38
$ws.Range("K53").Value = 38.24

# Row 54: "2021-01-09" / "17 Jan -- 23 Jan 2021" prediction
$ws.Range("A54").Value = "'2021-01-09"
$ws.Range("A54").Style = "Normal"
$ws.Range("B54").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D54").Value = 2211.75
$ws.Range("F54").Value = "KNN"

# Row 55: "2021-01-09" / "24 Jan -- 30 Jan 2021" prediction
$ws.Range("A55").Value = "'2021-01-09"
$ws.Range("A55").Style = "Normal"
$ws.Range("B55").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D55").Value = 2249.59
$ws.Range("F55").Value = "KNN"
